$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18807
$ws.Range("I21").Value = 8250
$ws.Range("J21").Value = 39921
$ws.Range("K21").Value = 8250
$ws.Range("L21").Value = 39921
$ws.Range("M21").Value = -7782
$ws.Range("N21").Value = -40857
$ws.Range("H23").Value = 18807
$ws.Range("I23").Value = 8250
$ws.Range("J23").Value = 39921
$ws.Range("K23").Value = 8250
$ws.Range("L23").Value = 39921
$ws.Range("M23").Value = -8016
$ws.Range("N23").Value = -40389
$ws.Range("H33").Value = 833
$ws.Range("I33").Value = 499
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 499
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -270
$ws.Range("N33").Value = -1458
$ws.Range("H86").Value = 6898.6
$ws.Range("I86").Value = 5831
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 5831
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -4708
$ws.Range("N86").Value = -10746
$ws.Range("H89").Value = 6898.6
$ws.Range("I89").Value = 5831
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 29155
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -23539
$ws.Range("N89").Value = -53732
$ws.Range("H97").Value = 2000
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H98").Value = 9797
$ws.Range("I98").Value = 6661.6665
$ws.Range("J98").Value = 14500
$ws.Range("K98").Value = 6661.6665
$ws.Range("L98").Value = 14500
$ws.Range("M98").Value = -5163.6665
$ws.Range("N98").Value = -17496
$ws.Range("H107").Value = 1885.8667
$ws.Range("I107").Value = 2146.3845
$ws.Range("K107").Value = 2146.3845
$ws.Range("M107").Value = -226.3845000000001
$ws.Range("H112").Value = 2824.1667
$ws.Range("J112").Value = 2824.1667
$ws.Range("L112").Value = 8472.500100000001
$ws.Range("N112").Value = -10688.5001
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 9797
$ws.Range("I122").Value = 6661.6665
$ws.Range("J122").Value = 14500
$ws.Range("K122").Value = 19984.9995
$ws.Range("L122").Value = 43500
$ws.Range("M122").Value = -17534.9995
$ws.Range("N122").Value = -48400
$ws.Range("H137").Value = 2852.8147
$ws.Range("I137").Value = 2693.7778
$ws.Range("J137").Value = 3170.889
$ws.Range("K137").Value = 8081.3334
$ws.Range("L137").Value = 9512.667000000001
$ws.Range("M137").Value = -5531.3334
$ws.Range("N137").Value = -14612.667
$ws.Range("H138").Value = 3153.6875
$ws.Range("I138").Value = 794.2
$ws.Range("J138").Value = 7086.1665
$ws.Range("K138").Value = 2382.6
$ws.Range("L138").Value = 21258.4995
$ws.Range("M138").Value = 2757.4
$ws.Range("N138").Value = -31538.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6669.6665
$ws.Range("I2").Value = 6669.6665
$ws.Range("K2").Value = 6669.6665
$ws.Range("M2").Value = -6556.6665
$ws.Range("H32").Value = 18691.777
$ws.Range("I32").Value = 18691.777
$ws.Range("K32").Value = 18691.777
$ws.Range("M32").Value = -18404.777
$ws.Range("H116").Value = 6669.6665
$ws.Range("I116").Value = 6669.6665
$ws.Range("K116").Value = 6669.6665
$ws.Range("M116").Value = -4375.6665
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6669.6665
$ws.Range("I3").Value = 6669.6665
$ws.Range("K3").Value = 6669.6665
$ws.Range("M3").Value = -6555.6665
$ws.Range("H64").Value = 8117
$ws.Range("J64").Value = 8117
$ws.Range("L64").Value = 8117
$ws.Range("N64").Value = -8567
$ws.Range("H67").Value = 8117
$ws.Range("J67").Value = 8117
$ws.Range("L67").Value = 8117
$ws.Range("N67").Value = -9677
$ws.Range("H105").Value = 9950
$ws.Range("I105").Value = 9950
$ws.Range("K105").Value = 9950
$ws.Range("M105").Value = -8203

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1907.7916
$ws.Range("I31").Value = 1799.4706
$ws.Range("J31").Value = 2170.8572
$ws.Range("K31").Value = 1799.4706
$ws.Range("L31").Value = 2170.8572
$ws.Range("M31").Value = -1504.4706
$ws.Range("N31").Value = -2760.8572
$ws.Range("H34").Value = 1907.7916
$ws.Range("I34").Value = 1799.4706
$ws.Range("J34").Value = 2170.8572
$ws.Range("K34").Value = 1799.4706
$ws.Range("L34").Value = 2170.8572
$ws.Range("M34").Value = -1597.4706
$ws.Range("N34").Value = -2574.8572
$ws.Range("H134").Value = 7649.3335
$ws.Range("I134").Value = 7179.4
$ws.Range("K134").Value = 21538.2
$ws.Range("M134").Value = -19003.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 949.5
$ws.Range("I5").Value = 949.5
$ws.Range("K5").Value = 2848.5
$ws.Range("M5").Value = -2736.5
$ws.Range("H17").Value = 1677.6154
$ws.Range("I17").Value = 146.66667
$ws.Range("J17").Value = 2136.9
$ws.Range("K17").Value = 440.00001
$ws.Range("L17").Value = 6410.700000000001
$ws.Range("M17").Value = -271.00001
$ws.Range("N17").Value = -6748.700000000001
$ws.Range("H68").Value = 14445
$ws.Range("I68").Value = 9002
$ws.Range("J68").Value = 19888
$ws.Range("K68").Value = 27006
$ws.Range("L68").Value = 59664
$ws.Range("M68").Value = -26195
$ws.Range("N68").Value = -61286
$ws.Range("H71").Value = 14445
$ws.Range("I71").Value = 9002
$ws.Range("J71").Value = 19888
$ws.Range("K71").Value = 81018
$ws.Range("L71").Value = 178992
$ws.Range("M71").Value = -76962
$ws.Range("N71").Value = -187104
$ws.Range("H123").Value = 4954.5454
$ws.Range("H124").Value = 5000
$ws.Range("J124").Value = 5000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820
$ws.Range("H129").Value = 2189.2727
$ws.Range("I129").Value = 1474.5
$ws.Range("K129").Value = 4423.5
$ws.Range("M129").Value = 576.5
$ws.Range("H131").Value = 2472.0476
$ws.Range("J131").Value = 2440.1667
$ws.Range("L131").Value = 7320.500100000001
$ws.Range("N131").Value = -17400.5001
$ws.Range("H135").Value = 949.5
$ws.Range("I135").Value = 949.5
$ws.Range("K135").Value = 8545.5
$ws.Range("M135").Value = -6010.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 10399.5
$ws.Range("J29").Value = 10399.5
$ws.Range("L29").Value = 10399.5
$ws.Range("N29").Value = -10979.5
$ws.Range("H80").Value = 9750
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 16500
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 16500
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -18496
$ws.Range("H83").Value = 9750
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 16500
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 82500
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -92484
$ws.Range("H126").Value = 4010.7778
$ws.Range("I126").Value = 4192.7144
$ws.Range("K126").Value = 12578.1432
$ws.Range("M126").Value = -10108.1432
$ws.Range("H132").Value = 3419
$ws.Range("I132").Value = 3336
$ws.Range("K132").Value = 10008
$ws.Range("M132").Value = -7478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6841
$ws.Range("I22").Value = 5895
$ws.Range("K22").Value = 5895
$ws.Range("M22").Value = -5600
$ws.Range("H27").Value = 6841
$ws.Range("I27").Value = 5895
$ws.Range("K27").Value = 5895
$ws.Range("M27").Value = -5788
$ws.Range("H82").Value = 3264.4
$ws.Range("I82").Value = 3264.4
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3264.4
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2903.4
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 3264.4
$ws.Range("I85").Value = 3264.4
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3264.4
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2016.4
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 44178.25
$ws.Range("I4").Value = 58952
$ws.Range("J4").Value = 19555.334
$ws.Range("K4").Value = 58952
$ws.Range("L4").Value = 19555.334
$ws.Range("M4").Value = -58839
$ws.Range("N4").Value = -19781.334
$ws.Range("H136").Value = 2019.24
$ws.Range("I136").Value = 2225.4
$ws.Range("K136").Value = 6676.200000000001
$ws.Range("M136").Value = -4126.200000000001